# Add the "NY Quits - Smokers Quit Line" rolodex entries (AOD + Hotline rows)
# to the Rolodex worksheet, right after the existing last row (119).
#
# Note: cells are written in this specific order (Website, Name, Phone,
# Hours, Type, Loc, Info) so that newly-minted shared-string entries land
# at the same indices the source workbook used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$name    = "NY Quits - Smokers Quit Line"
$phone   = "te:1-866-697-8487"
$website = "https://www.nysmokefree.com/"
$hours   = "M-Th: 9AM-9PM & Fri-Sun: 1AM-5PM"
$info    = "Interested in quitting, or need support managing nicotine cravings? Call the NYS Smokers' Quitlin to apply for a free starter kit of nicotine medications and to talk to a quit coach. FREE!"

# Row 120 - AOD entry
$ws.Range("D120").Value = $website
$ws.Range("A120").Value = $name
$ws.Range("C120").Value = $phone
$ws.Range("E120").Value = $hours
$ws.Range("G120").Value = "AOD"
$ws.Range("G120").NumberFormat = "@"
$ws.Range("H120").Value = "Off"
$ws.Range("I120").Value = $info

# Row 121 - Hotline entry
$ws.Range("A121").Value = $name
$ws.Range("C121").Value = $phone
$ws.Range("D121").Value = $website
$ws.Range("E121").Value = $hours
$ws.Range("G121").Value = "Hotline"
$ws.Range("G121").NumberFormat = "@"
$ws.Range("H121").Value = "Off"
$ws.Range("I121").Value = $info

# Restore the selection to match the post-edit cursor position
$ws.Range("F129").Select()
